# Prepend 8 new rows of sensor readings (timestamps 0-700) ahead of the
# existing "falling" dataset, and append 2 more rows at the end
# (timestamps 2800-2900) -- "the changes we made on may 9th".
# The full resulting table (rows 2-31, A:timestamp .. H:gz) is written
# out explicitly below since every timestamp/row shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: timestamp, ax, ay, az, gx, gy, gz  (label is always "falling")
$data = @(
    @(0, -0.1574550628662107, 0.1197370529174804, -0.3381910085678101, -0.0337503030896186, 0.09865473955869671, -0.0665843114256858),
    @(100, -0.3228707909584047, 0.1800707578659058, -0.2396511320024727, -0.0201585534960031, -0.0109955742955207, -0.0201585534960031),
    @(200, -0.3687732219696044, 0.22406005859375, -0.1331830434501169, 0.0261144898831844, -0.0386372283101081, 0.0106901414692401),
    @(300, -0.3819821834564209, 0.1978947639465331, -0.2025355234742165, 0.0059559359215199, -0.0545197241008281, 0.0339030213654041),
    @(400, -0.3319910526275633, 0.2015152931213379, -0.2596112400293351, 0.0154243474826216, -0.0274889357388019, -0.0103847095742821),
    @(500, -0.3998157501220704, 0.1588554024696349, -0.4224489629268647, -0.0024434609804302, -0.0116064399480819, -0.0108428578823804),
    @(600, -0.4262521266937255, -0.03637795448303255, -0.570676848292351, 0.0427605658769607, 0.0514653958380222, -0.0694859251379966),
    @(700, -0.5638278126716616, -0.07290470600128157, -1.039905160665513, 0.0119118718430399, -0.0441350154578685, 0.0128281703218817),
    @(800, -0.4378628730773919, -0.1533573150634768, -1.335271286964417, 0.0447458773851394, -0.1001818999648094, 0.0022907445672899),
    @(900, -0.2312893986701963, -0.572000253200532, -1.991752985119822, 0.0861319974064827, -0.4196644127368927, 0.0219911485910415),
    @(1000, -0.05911626815795876, -1.518646645545961, -2.819378316402437, 0.1461495161056518, -0.6041457056999207, 0.3318525552749634),
    @(1100, -0.867481225728992, -2.808535265922549, -3.745611393451694, 0.2295326143503189, -0.7486153841018677, 0.0561996027827262),
    @(1200, -2.066653013229371, -4.002454280853272, -4.322827577590942, 0.7399105429649353, 0.0711658000946044, -1.335504174232483),
    @(1300, -1.651154989004131, -3.045457839965813, -5.523038685321813, 4.26658821105957, 2.570826292037964, -3.828750610351562),
    @(1400, -3.538045167922988, -0.1250472068786557, -6.249881243705746, 0.669966459274292, -0.3785837292671203, -2.870913982391357),
    @(1500, 21.8330995559694, -0.8298350632190838, 8.03648120164879, -1.501964926719666, 0.1320996135473251, -1.118036150932312),
    @(1600, 36.38044357299789, -2.828326940536493, 15.55662982463829, 0.0074830991216003, 0.1238529309630394, 0.8869763612747192),
    @(1700, -0.2419787645339966, -1.207208663225174, -0.5454361662268639, 0.0334448739886283, -0.0073303831741213, -0.1014036312699318),
    @(1800, 0.3278889179229798, -0.6977443240582908, -0.3904904812574392, -0.4305072724819183, -0.0621555373072624, 0.0183259565383195),
    @(1900, 0.6186044454574491, -0.1085777457803506, -0.6750323086977017, -1.069319605827332, 0.0797179117798805, -0.1148426681756973),
    @(2000, -1.584584712982176, -0.268941026926038, -0.4205022692680305, 0.2054034322500229, -0.0172569435089826, -0.1675297915935516),
    @(2100, -1.204359865188596, 0.3030472241342062, 0.2759267613291676, 0.0361937656998634, 0.0485637858510017, 0.1411098688840866),
    @(2200, -0.6684749126434326, 0.3942791000008581, -1.052078396081925, -0.1175915598869323, 0.0036651915870606, 0.0314595587551593),
    @(2300, -0.5069202423095689, 0.8594004645943611, -0.3887748748064006, -0.090866208076477, 0.0201585534960031, 0.0360410511493682),
    @(2400, -0.1816680908203125, 0.3003633975982678, 0.3908397436141963, -0.007177666760981, -0.072540245950222, -0.0835358202457428),
    @(2500, -0.2214233398437506, 0.5106331750750528, 0.1910619698464848, 0.0241291765123605, -0.0126754539087414, -0.0610865242779254),
    @(2600, -0.3159080505371084, 0.2065896064043048, -0.2784154430031768, 0.0061086523346602, -0.0195476878434419, -0.0714712366461753),
    @(2700, -0.06431400775909421, 0.2912909723818302, -0.08320332318544384, 0.0032070425804704, 0.0010690141934901, -0.0012217304902151),
    @(2800, -0.003773307800292898, 0.2548617899417878, -0.1365976125001907, 0.0019853119738399, -0.0036651915870606, 0.0067195175215601),
    @(2900, 0.01674773693084754, 0.2881045341491697, -0.1470022119581697, 0.0003054326225537, -0.00167987938039, 0.0142026171088218)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "falling"
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, 3 + $c).Value = $row[1 + $c]
    }
}
